$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Brasil'
$ws.Range("B2").Value = 0.5306999999999999
$ws.Range("C2").Value = 0.2861
$ws.Range("D2").Value = 0.1527
$ws.Range("E2").Value = 0.0305
$ws.Range("F2").Value = 0.8168
$ws.Range("G2").Value = 0.5181
$ws.Range("H2").Value = 0.3214
$ws.Range("I2").Value = 0.2048
$ws.Range("J2").Value = 0.125

$ws.Range("A3").Value = 'Argentina'
$ws.Range("B3").Value = 0.5683
$ws.Range("C3").Value = 0.2453
$ws.Range("D3").Value = 0.1315
$ws.Range("E3").Value = 0.0549
$ws.Range("F3").Value = 0.8136
$ws.Range("G3").Value = 0.4951
$ws.Range("H3").Value = 0.3026
$ws.Range("I3").Value = 0.1729
$ws.Range("J3").Value = 0.0985

$ws.Range("A4").Value = 'Holanda'
$ws.Range("B4").Value = 0.5698
$ws.Range("C4").Value = 0.2812
$ws.Range("D4").Value = 0.1136
$ws.Range("E4").Value = 0.0354
$ws.Range("F4").Value = 0.851
$ws.Range("G4").Value = 0.5405
$ws.Range("H4").Value = 0.3107
$ws.Range("I4").Value = 0.1692
$ws.Range("J4").Value = 0.0945

$ws.Range("A5").Value = 'Espanha'
$ws.Range("B5").Value = 0.4417
$ws.Range("C5").Value = 0.2961
$ws.Range("D5").Value = 0.1726
$ws.Range("E5").Value = 0.0896
$ws.Range("F5").Value = 0.7378
$ws.Range("G5").Value = 0.4371
$ws.Range("H5").Value = 0.2334
$ws.Range("I5").Value = 0.1311
$ws.Range("J5").Value = 0.0741

$ws.Range("A6").Value = 'França'
$ws.Range("B6").Value = 0.4347
$ws.Range("C6").Value = 0.3115
$ws.Range("D6").Value = 0.1658
$ws.Range("E6").Value = 0.08799999999999999
$ws.Range("F6").Value = 0.7462
$ws.Range("G6").Value = 0.4374
$ws.Range("H6").Value = 0.2521
$ws.Range("I6").Value = 0.1311
$ws.Range("J6").Value = 0.0674

$ws.Range("A7").Value = 'Portugal'
$ws.Range("B7").Value = 0.4526
$ws.Range("C7").Value = 0.3261
$ws.Range("D7").Value = 0.1858
$ws.Range("E7").Value = 0.0355
$ws.Range("F7").Value = 0.7786999999999999
$ws.Range("G7").Value = 0.4056
$ws.Range("H7").Value = 0.2295
$ws.Range("I7").Value = 0.1277
$ws.Range("J7").Value = 0.065

$ws.Range("A8").Value = 'Bélgica'
$ws.Range("B8").Value = 0.4235
$ws.Range("C8").Value = 0.2895
$ws.Range("D8").Value = 0.1785
$ws.Range("E8").Value = 0.1085
$ws.Range("F8").Value = 0.713
$ws.Range("G8").Value = 0.3998
$ws.Range("H8").Value = 0.2144
$ws.Range("I8").Value = 0.1229
$ws.Range("J8").Value = 0.0634

$ws.Range("A9").Value = 'Dinamarca'
$ws.Range("B9").Value = 0.3754
$ws.Range("C9").Value = 0.326
$ws.Range("D9").Value = 0.1947
$ws.Range("E9").Value = 0.1039
$ws.Range("F9").Value = 0.7014
$ws.Range("G9").Value = 0.3934
$ws.Range("H9").Value = 0.2166
$ws.Range("I9").Value = 0.1039
$ws.Range("J9").Value = 0.0501

$ws.Range("A10").Value = 'Inglaterra'
$ws.Range("B10").Value = 0.4416
$ws.Range("C10").Value = 0.2718
$ws.Range("D10").Value = 0.1777
$ws.Range("E10").Value = 0.1089
$ws.Range("F10").Value = 0.7134
$ws.Range("G10").Value = 0.3987
$ws.Range("H10").Value = 0.2126
$ws.Range("I10").Value = 0.102
$ws.Range("J10").Value = 0.0497

$ws.Range("A11").Value = 'Alemanha'
$ws.Range("B11").Value = 0.3386
$ws.Range("C11").Value = 0.3167
$ws.Range("D11").Value = 0.217
$ws.Range("E11").Value = 0.1277
$ws.Range("F11").Value = 0.6553
$ws.Range("G11").Value = 0.3565
$ws.Range("H11").Value = 0.1804
$ws.Range("I11").Value = 0.0964
$ws.Range("J11").Value = 0.0472

$ws.Range("A12").Value = 'Uruguai'
$ws.Range("B12").Value = 0.3639
$ws.Range("C12").Value = 0.3489
$ws.Range("D12").Value = 0.2409
$ws.Range("E12").Value = 0.0463
$ws.Range("F12").Value = 0.7128
$ws.Range("G12").Value = 0.3398
$ws.Range("H12").Value = 0.1757
$ws.Range("I12").Value = 0.08740000000000001
$ws.Range("J12").Value = 0.0433

$ws.Range("A13").Value = 'Croácia'
$ws.Range("B13").Value = 0.3187
$ws.Range("C13").Value = 0.3004
$ws.Range("D13").Value = 0.2263
$ws.Range("E13").Value = 0.1546
$ws.Range("F13").Value = 0.6191
$ws.Range("G13").Value = 0.3148
$ws.Range("H13").Value = 0.1429
$ws.Range("I13").Value = 0.0732
$ws.Range("J13").Value = 0.0348

$ws.Range("A14").Value = 'Equador'
$ws.Range("B14").Value = 0.3097
$ws.Range("C14").Value = 0.3825
$ws.Range("D14").Value = 0.2098
$ws.Range("E14").Value = 0.098
$ws.Range("F14").Value = 0.6922
$ws.Range("G14").Value = 0.3697
$ws.Range("H14").Value = 0.1684
$ws.Range("I14").Value = 0.0752
$ws.Range("J14").Value = 0.033

$ws.Range("A15").Value = 'Suíça'
$ws.Range("B15").Value = 0.2245
$ws.Range("C15").Value = 0.3262
$ws.Range("D15").Value = 0.3406
$ws.Range("E15").Value = 0.1087
$ws.Range("F15").Value = 0.5507
$ws.Range("G15").Value = 0.2736
$ws.Range("H15").Value = 0.1318
$ws.Range("I15").Value = 0.0655
$ws.Range("J15").Value = 0.0296

$ws.Range("A16").Value = 'Sérvia'
$ws.Range("B16").Value = 0.2319
$ws.Range("C16").Value = 0.3306
$ws.Range("D16").Value = 0.3299
$ws.Range("E16").Value = 0.1076
$ws.Range("F16").Value = 0.5625
$ws.Range("G16").Value = 0.2733
$ws.Range("H16").Value = 0.129
$ws.Range("I16").Value = 0.0615
$ws.Range("J16").Value = 0.0265

$ws.Range("A17").Value = 'México'
$ws.Range("B17").Value = 0.1832
$ws.Range("C17").Value = 0.2925
$ws.Range("D17").Value = 0.2998
$ws.Range("E17").Value = 0.2245
$ws.Range("F17").Value = 0.4757
$ws.Range("G17").Value = 0.2014
$ws.Range("H17").Value = 0.0915
$ws.Range("I17").Value = 0.0341
$ws.Range("J17").Value = 0.0142

$ws.Range("A18").Value = 'Estados Unidos'
$ws.Range("B18").Value = 0.21
$ws.Range("C18").Value = 0.2537
$ws.Range("D18").Value = 0.267
$ws.Range("E18").Value = 0.2693
$ws.Range("F18").Value = 0.4637
$ws.Range("G18").Value = 0.2086
$ws.Range("H18").Value = 0.0872
$ws.Range("I18").Value = 0.0324
$ws.Range("J18").Value = 0.0136

$ws.Range("A19").Value = 'Polônia'
$ws.Range("B19").Value = 0.1841
$ws.Range("C19").Value = 0.3038
$ws.Range("D19").Value = 0.2967
$ws.Range("E19").Value = 0.2154
$ws.Range("F19").Value = 0.4879
$ws.Range("G19").Value = 0.2127
$ws.Range("H19").Value = 0.0916
$ws.Range("I19").Value = 0.0351
$ws.Range("J19").Value = 0.0125

$ws.Range("A20").Value = 'Coreia do Sul'
$ws.Range("B20").Value = 0.1751
$ws.Range("C20").Value = 0.2866
$ws.Range("D20").Value = 0.4237
$ws.Range("E20").Value = 0.1146
$ws.Range("F20").Value = 0.4617
$ws.Range("G20").Value = 0.168
$ws.Range("H20").Value = 0.0692
$ws.Range("I20").Value = 0.0286
$ws.Range("J20").Value = 0.0098

$ws.Range("A21").Value = 'País de Gales'
$ws.Range("B21").Value = 0.1891
$ws.Range("C21").Value = 0.2529
$ws.Range("D21").Value = 0.2733
$ws.Range("E21").Value = 0.2847
$ws.Range("F21").Value = 0.442
$ws.Range("G21").Value = 0.1881
$ws.Range("H21").Value = 0.0733
$ws.Range("I21").Value = 0.0229
$ws.Range("J21").Value = 0.0083

$ws.Range("A22").Value = 'Irã'
$ws.Range("B22").Value = 0.1593
$ws.Range("C22").Value = 0.2216
$ws.Range("D22").Value = 0.282
$ws.Range("E22").Value = 0.3371
$ws.Range("F22").Value = 0.3809
$ws.Range("G22").Value = 0.1542
$ws.Range("H22").Value = 0.062
$ws.Range("I22").Value = 0.0202
$ws.Range("J22").Value = 0.0074

$ws.Range("A23").Value = 'Japão'
$ws.Range("B23").Value = 0.1283
$ws.Range("C23").Value = 0.2153
$ws.Range("D23").Value = 0.3112
$ws.Range("E23").Value = 0.3452
$ws.Range("F23").Value = 0.3436
$ws.Range("G23").Value = 0.1393
$ws.Range("H23").Value = 0.0466
$ws.Range("I23").Value = 0.019
$ws.Range("J23").Value = 0.0062

$ws.Range("A24").Value = 'Marrocos'
$ws.Range("B24").Value = 0.1237
$ws.Range("C24").Value = 0.195
$ws.Range("D24").Value = 0.2981
$ws.Range("E24").Value = 0.3832
$ws.Range("F24").Value = 0.3187
$ws.Range("G24").Value = 0.1212
$ws.Range("H24").Value = 0.0419
$ws.Range("I24").Value = 0.0159
$ws.Range("J24").Value = 0.0057

$ws.Range("A25").Value = 'Canadá'
$ws.Range("B25").Value = 0.1341
$ws.Range("C25").Value = 0.2151
$ws.Range("D25").Value = 0.2971
$ws.Range("E25").Value = 0.3537
$ws.Range("F25").Value = 0.3492
$ws.Range("G25").Value = 0.1359
$ws.Range("H25").Value = 0.0496
$ws.Range("I25").Value = 0.0179
$ws.Range("J25").Value = 0.0057

$ws.Range("A26").Value = 'Austrália'
$ws.Range("B26").Value = 0.0955
$ws.Range("C26").Value = 0.18
$ws.Range("D26").Value = 0.3107
$ws.Range("E26").Value = 0.4138
$ws.Range("F26").Value = 0.2755
$ws.Range("G26").Value = 0.09660000000000001
$ws.Range("H26").Value = 0.0341
$ws.Range("I26").Value = 0.0105
$ws.Range("J26").Value = 0.0046

$ws.Range("A27").Value = 'Costa Rica'
$ws.Range("B27").Value = 0.0914
$ws.Range("C27").Value = 0.1719
$ws.Range("D27").Value = 0.2992
$ws.Range("E27").Value = 0.4375
$ws.Range("F27").Value = 0.2633
$ws.Range("G27").Value = 0.0954
$ws.Range("H27").Value = 0.0306
$ws.Range("I27").Value = 0.0115
$ws.Range("J27").Value = 0.0034

$ws.Range("A28").Value = 'Tunísia'
$ws.Range("B28").Value = 0.0944
$ws.Range("C28").Value = 0.1825
$ws.Range("D28").Value = 0.3288
$ws.Range("E28").Value = 0.3943
$ws.Range("F28").Value = 0.2769
$ws.Range("G28").Value = 0.0969
$ws.Range("H28").Value = 0.0355
$ws.Range("I28").Value = 0.011
$ws.Range("J28").Value = 0.0033

$ws.Range("A29").Value = 'Arábia Saudita'
$ws.Range("B29").Value = 0.0644
$ws.Range("C29").Value = 0.1584
$ws.Range("D29").Value = 0.272
$ws.Range("E29").Value = 0.5052
$ws.Range("F29").Value = 0.2228
$ws.Range("G29").Value = 0.0665
$ws.Range("H29").Value = 0.0211
$ws.Range("I29").Value = 0.0059
$ws.Range("J29").Value = 0.0018

$ws.Range("A30").Value = 'Senegal'
$ws.Range("B30").Value = 0.0707
$ws.Range("C30").Value = 0.1969
$ws.Range("D30").Value = 0.3617
$ws.Range("E30").Value = 0.3707
$ws.Range("F30").Value = 0.2676
$ws.Range("G30").Value = 0.091
$ws.Range("H30").Value = 0.0277
$ws.Range("I30").Value = 0.0069
$ws.Range("J30").Value = 0.0013

$ws.Range("A31").Value = 'Catar'
$ws.Range("B31").Value = 0.0498
$ws.Range("C31").Value = 0.1394
$ws.Range("D31").Value = 0.3149
$ws.Range("E31").Value = 0.4959
$ws.Range("F31").Value = 0.1892
$ws.Range("G31").Value = 0.0492
$ws.Range("H31").Value = 0.013
$ws.Range("I31").Value = 0.0028
$ws.Range("J31").Value = 0.0001

$ws.Range("A32").Value = 'Camarões'
$ws.Range("B32").Value = 0.0129
$ws.Range("C32").Value = 0.0571
$ws.Range("D32").Value = 0.1768
$ws.Range("E32").Value = 0.7532
$ws.Range("F32").Value = 0.07000000000000001
$ws.Range("G32").Value = 0.0152
$ws.Range("H32").Value = 0.0028
$ws.Range("I32").Value = 0.0005
$ws.Range("J32").Value = 0

$ws.Range("A33").Value = 'Gana'
$ws.Range("B33").Value = 0.008399999999999999
$ws.Range("C33").Value = 0.0384
$ws.Range("D33").Value = 0.1496
$ws.Range("E33").Value = 0.8036
$ws.Range("F33").Value = 0.0468
$ws.Range("G33").Value = 0.0064
$ws.Range("H33").Value = 0.0008
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
